# application-phase.xlsx — add a "Version" column.
#
# Before:
#   A: Code   B: Description              C: Definition
#   S1 / M1 / M2 / O1 rows with their description + definition text.
#
# After (per commit "Updated with version numbers, updated JSON structure,
# refined codesets"): a new "Version" column is inserted as the new first
# column (A), pushing Code/Description/Definition one column to the right
# (B/C/D). Every data row gets the version value "1.0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: Version, Code, Description, Definition
$ws.Range("A1").Value = "Version"
$ws.Range("B1").Value = "Code"
$ws.Range("C1").Value = "Description"
$ws.Range("D1").Value = "Definition"

# The version numbers ("1.0") must stay plain text (matching the source
# workbook's "1.0" shared string) instead of Excel's default behaviour of
# silently parsing a numeric-looking string into the number 1. Mark the
# column as text before writing, then drop the formatting override again
# so the cells end up back on the default (unstyled) cell format.
$ws.Range("A2:A5").NumberFormat = "@"

# Row 2 - S1 / Single Stage Application
$ws.Range("A2").Value = "1.0"
$ws.Range("B2").Value = "S1"
$ws.Range("C2").Value = "Single Stage Application"
$ws.Range("D2").Value = "A full application for a single stage competitive, negotiated, or on demand round"

# Row 3 - M1 / Initial Application
$ws.Range("A3").Value = "1.0"
$ws.Range("B3").Value = "M1"
$ws.Range("C3").Value = "Initial Application"
$ws.Range("D3").Value = "An initial expression of interest in a multi-stage competitive round"

# Row 4 - M2 / Detailed Application
$ws.Range("A4").Value = "1.0"
$ws.Range("B4").Value = "M2"
$ws.Range("C4").Value = "Detailed Application"
$ws.Range("D4").Value = "An detailed application subsequent to an accepted expression of interest in a multi-stage competitive round"

# Row 5 - O1 / Other Application Type
$ws.Range("A5").Value = "1.0"
$ws.Range("B5").Value = "O1"
$ws.Range("C5").Value = "Other Application Type"
$ws.Range("D5").Value = "A form of application not otherwise specified in this code set"

# Drop the temporary text-format override again (back to the default style).
$ws.Range("A2:A5").Style = "Normal"
